$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 151
$ws.Range("F3").Value = 162
$ws.Range("F4").Value = 2144
$ws.Range("F5").Value = 4301
$ws.Range("F6").Value = 566
$ws.Range("F7").Value = 1063
$ws.Range("F8").Value = 1327
$ws.Range("F11").Value = 2227
$ws.Range("F12").Value = 401
$ws.Range("F13").Value = 666087
$ws.Range("F14").Value = 1659
$ws.Range("F15").Value = 555
$ws.Range("F16").Value = 1479
$ws.Range("F19").Value = 1292
$ws.Range("F20").Value = 2279
$ws.Range("F21").Value = 1158
$ws.Range("F22").Value = 2712
$ws.Range("F23").Value = 1570
$ws.Range("F24").Value = 854
$ws.Range("F25").Value = 1564
$ws.Range("F27").Value = 1094
$ws.Range("F28").Value = 347
$ws.Range("F29").Value = 1100
$ws.Range("F30").Value = 45
$ws.Range("F31").Value = 86
$ws.Range("F32").Value = 2040
$ws.Range("F33").Value = 1402
$ws.Range("F34").Value = 586
$ws.Range("F35").Value = 1308
$ws.Range("F36").Value = 2749
$ws.Range("F37").Value = 9
$ws.Range("F38").Value = 1151
$ws.Range("F41").Value = 2614
$ws.Range("F42").Value = 214
$ws.Range("F43").Value = 1001
$ws.Range("F44").Value = 3165
$ws.Range("F45").Value = 1020
$ws.Range("F48").Value = 161
$ws.Range("F50").Value = 19

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 69
$ws.Range("F10").Value = 489
$ws.Range("F11").Value = 144795
$ws.Range("F12").Value = 144795
$ws.Range("F13").Value = 11
$ws.Range("F15").Value = 22
$ws.Range("F17").Value = 97
$ws.Range("F18").Value = 233
$ws.Range("F19").Value = 338
$ws.Range("F21").Value = 419
$ws.Range("F22").Value = 162
$ws.Range("F24").Value = 95
$ws.Range("F26").Value = 606
$ws.Range("F28").Value = 1
$ws.Range("F31").Value = 363
$ws.Range("F32").Value = 278
$ws.Range("F34").Value = 64
$ws.Range("F35").Value = 64
$ws.Range("F38").Value = 216
$ws.Range("F43").Value = 9

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F5").Value = 251
$ws.Range("F7").Value = 829
$ws.Range("F8").Value = 1207
$ws.Range("F9").Value = 642
$ws.Range("F10").Value = 1610
$ws.Range("F11").Value = 477
$ws.Range("F12").Value = 125
$ws.Range("F13").Value = 1978

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 829
$ws.Range("F3").Value = 1207
$ws.Range("F4").Value = 642
$ws.Range("F5").Value = 151
$ws.Range("F6").Value = 1610
$ws.Range("F7").Value = 162
$ws.Range("F8").Value = 2144
$ws.Range("F9").Value = 125
$ws.Range("F10").Value = 1978
$ws.Range("F11").Value = 4301
$ws.Range("F12").Value = 566
$ws.Range("F13").Value = 1327
$ws.Range("F16").Value = 2227
$ws.Range("F17").Value = 401
$ws.Range("F18").Value = 666087
$ws.Range("F20").Value = 489
$ws.Range("F21").Value = 1659
$ws.Range("F22").Value = 144795
$ws.Range("F23").Value = 1479
$ws.Range("F26").Value = 1292
$ws.Range("F27").Value = 2279
$ws.Range("F28").Value = 1158
$ws.Range("F29").Value = 2712
$ws.Range("F30").Value = 1570
$ws.Range("F31").Value = 854
$ws.Range("F32").Value = 22
$ws.Range("F33").Value = 1564
$ws.Range("F35").Value = 163
$ws.Range("F36").Value = 1094
$ws.Range("F37").Value = 1100
$ws.Range("F38").Value = 1402
$ws.Range("F39").Value = 1308
$ws.Range("F40").Value = 2749
$ws.Range("F41").Value = 9
$ws.Range("F42").Value = 1151
$ws.Range("F43").Value = 363
$ws.Range("F44").Value = 278
$ws.Range("F45").Value = 64
$ws.Range("F46").Value = 2614
$ws.Range("F47").Value = 1001
$ws.Range("F48").Value = 3165
$ws.Range("F49").Value = 216
$ws.Range("F50").Value = 1020
$ws.Range("F51").Value = 161
$ws.Range("F53").Value = 19
